# Actualización desde MV -datos-
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the last existing data point (row 272, 28-09-2021)
$ws.Cells.Item(272, 2).Value = 189.2581

# Append the new daily observations (rows 273-279)
$dates  = @("29-09-2021", "30-09-2021", "01-10-2021", "02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021")
$values = @(193.4353, 195.1864, 194.3698, 193.2917, 185.9804, 184.7754, 183.0694)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 273 + $i

    # Write the date through a formula and paste back as a value so the
    # text ("dd-mm-yyyy") lands as a plain shared string instead of being
    # auto-recognised as a date serial number, and without touching the
    # cell's style (no NumberFormat / quote-prefix residue).
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = '="' + $dates[$i] + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null

    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$excel.CutCopyMode = 0
